# Update the hyperparameter-search results sheet for the newly-finished
# training run (target_col == previous_concussions).
#  - Row 2: refreshed selector object id + model hyperparams/scores (in place)
#  - Old row 3 (max_depth=2, min_samples_split=3, n_estimators=100) is gone;
#    the row that used to be row 4 now occupies row 3 with updated
#    hyperparams/scores, and a brand-new row 4 (max_depth=3) is added.
#  - Rows 5 and 6 get refreshed hyperparams/scores in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$B2 = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98e0eeb3a0>),
                ('model',
                 RandomForestClassifier(class_weight='balanced', max_depth=1,
                                        max_features='log2', min_samples_leaf=3,
                                        min_samples_split=5, n_estimators=350,
                                        random_state=42))])
"@
$ws.Range("B2").Value = $B2

$C2 = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f99170fc070>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 350, 'model__max_depth': 1, 'model__min_samples_split': 5, 'model__min_samples_leaf': 3, 'model__max_features': 'log2', 'model__class_weight': 'balanced'}
"@
$ws.Range("C2").Value = $C2

$ws.Range("D2").Value = 0.5701098397198141

$ws.Range("G2").Value = 0.6900503070862122

$B3 = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98e11b24f0>),
                ('model',
                 RandomForestClassifier(class_weight='balanced', max_depth=2,
                                        min_samples_leaf=11,
                                        min_samples_split=4, n_estimators=250,
                                        random_state=42))])
"@
$ws.Range("B3").Value = $B3

$C3 = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f98e0f28970>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 250, 'model__max_depth': 2, 'model__min_samples_split': 4, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__class_weight': 'balanced'}
"@
$ws.Range("C3").Value = $C3

$ws.Range("D3").Value = 0.5575734726678326

$ws.Range("G3").Value = 0.7117255504352279

$ws.Range("H3").Value = 0.4207459207459207

$J3 = @"
[1 1 1 1 1 1 1 1 0 1 0 0 0 0 1 1 0 0 1 0 0 0 0 0]
"@
$ws.Range("J3").Value = $J3

$B4 = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'D...
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f98e1010b20>),
                ('model',
                 RandomForestClassifier(class_weight='balanced', max_depth=3,
                                        max_features='log2', min_samples_leaf=8,
                                        min_samples_split=10, n_estimators=200,
                                        random_state=42))])
"@
$ws.Range("B4").Value = $B4

$C4 = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f98642c3700>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 200, 'model__max_depth': 3, 'model__min_samples_split': 10, 'model__min_samples_leaf': 8, 'model__max_features': 'log2', 'model__class_weight': 'balanced'}
"@
$ws.Range("C4").Value = $C4

$ws.Range("D4").Value = 0.5377129631410155

$ws.Range("G4").Value = 0.8044708545557442

$B5 = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 RandomForestClassifier(class_weight='balanced', max_depth=1,
                                        min_samples_leaf=2, min_samples_split=8,
                                        n_estimators=450, random_state=42))])
"@
$ws.Range("B5").Value = $B5

$C5 = @"
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 450, 'model__max_depth': 1, 'model__min_samples_split': 8, 'model__min_samples_leaf': 2, 'model__max_features': 'sqrt', 'model__class_weight': 'balanced'}
"@
$ws.Range("C5").Value = $C5

$ws.Range("D5").Value = 0.5971479110919143

$ws.Range("G5").Value = 0.7435179549604384

$ws.Range("H5").Value = 0.5440917107583775

$J5 = @"
[0 0 1 0 0 1 1 1 0 1 0 1 1 0 1 1 0 1 1 0 1 1 0 0]
"@
$ws.Range("J5").Value = $J5

$B6 = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 RandomForestClassifier(class_weight='balanced', max_depth=1,
                                        min_samples_leaf=5, min_samples_split=6,
                                        n_estimators=300, random_state=42))])
"@
$ws.Range("B6").Value = $B6

$C6 = @"
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 300, 'model__max_depth': 1, 'model__min_samples_split': 6, 'model__min_samples_leaf': 5, 'model__max_features': 'sqrt', 'model__class_weight': 'balanced'}
"@
$ws.Range("C6").Value = $C6

$ws.Range("D6").Value = 0.5846973785700398

$ws.Range("G6").Value = 0.7531596350773623

$ws.Range("A2:A6").EntireRow.AutoFit()
